$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.229.82'

$ws.Cells.Item(3, 4).Value = '2.643.83'
$ws.Cells.Item(3, 5).Value = '  -0.07%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).Value = '596.21'
$ws.Cells.Item(5, 5).Value = '  -0.32%  '

$ws.Cells.Item(6, 4).Value = '158.89'
$ws.Cells.Item(6, 5).Value = '  +2.59%  '

$ws.Cells.Item(7, 5).Value = '  +0.01%  '

$ws.Cells.Item(8, 4).Value = '0.542'
$ws.Cells.Item(8, 5).Value = '  -0.90%  '

$ws.Cells.Item(9, 5).Value = '  -2.48%  '

$ws.Cells.Item(10, 5).Value = '  -1.08%  '

$ws.Cells.Item(11, 4).Value = '5.28'
$ws.Cells.Item(11, 5).Value = '  +0.23%  '

$ws.Cells.Item(12, 5).Value = '  -0.91%  '

$ws.Cells.Item(13, 4).Value = '27.92'
$ws.Cells.Item(13, 5).Value = '  -1.12%  '

$ws.Cells.Item(14, 4).Value = '3.126.79'
$ws.Cells.Item(14, 5).Value = '  +0.24%  '

$ws.Cells.Item(15, 5).Value = '  -3.07%  '

$ws.Cells.Item(16, 4).Value = '68.084.81'
$ws.Cells.Item(16, 5).Value = '  -0.34%  '

$ws.Cells.Item(17, 4).Value = '2.617.72'
$ws.Cells.Item(17, 5).Value = '  -0.80%  '

$ws.Cells.Item(18, 4).Value = '11.35'
$ws.Cells.Item(18, 5).Value = '  -0.81%  '

$ws.Cells.Item(19, 4).Value = '358.42'
$ws.Cells.Item(19, 5).Value = '  -1.78%  '

$ws.Cells.Item(20, 4).Value = '7.38'
$ws.Cells.Item(20, 5).Value = '  -1.88%  '

$ws.Cells.Item(21, 4).Value = '4.41'
$ws.Cells.Item(21, 5).Value = '  +0.26%  '

$ws.Cells.Item(22, 4).Value = '4.74'
$ws.Cells.Item(22, 5).Value = '  -3.35%  '

$ws.Cells.Item(23, 5).Value = '  -0.65%  '

$ws.Cells.Item(24, 4).Value = '74.59'
$ws.Cells.Item(24, 5).Value = '  +0.18%  '

$ws.Cells.Item(25, 5).Value = '  -0.07%  '

$ws.Cells.Item(26, 4).Value = '9.71'
$ws.Cells.Item(26, 5).Value = '  -1.04%  '

$ws.Cells.Item(27, 4).Value = '2.779.51'
$ws.Cells.Item(27, 5).Value = '  +0.15%  '

$ws.Cells.Item(28, 5).Value = '  -3.03%  '

$ws.Cells.Item(29, 5).Value = '  +0.20%  '

$ws.Cells.Item(30, 4).Value = '558.62'
$ws.Cells.Item(30, 5).Value = '  -2.58%  '

$ws.Cells.Item(31, 4).Value = '7.99'
$ws.Cells.Item(31, 5).Value = '  -2.63%  '

$ws.Cells.Item(32, 4).Value = '1.39'
$ws.Cells.Item(32, 5).Value = '  -2.48%  '

$ws.Cells.Item(33, 4).Value = '1.88'
$ws.Cells.Item(33, 5).Value = '  +0.77%  '

$ws.Cells.Item(34, 5).Value = '  +1.69%  '

$ws.Cells.Item(35, 5).Value = '  +0.04%  '

$ws.Cells.Item(36, 5).Value = '  -2.99%  '

$ws.Cells.Item(37, 4).Value = '159.75'
$ws.Cells.Item(37, 5).Value = '  -0.69%  '

$ws.Cells.Item(39, 4).Value = '0.370'
$ws.Cells.Item(39, 5).Value = '  -1.18%  '

$ws.Cells.Item(40, 4).Value = '1.86'
$ws.Cells.Item(40, 5).Value = '  -2.69%  '

$ws.Cells.Item(41, 4).Value = '5.32'
$ws.Cells.Item(41, 5).Value = '  -2.05%  '

$ws.Cells.Item(42, 4).Value = '2.60'
$ws.Cells.Item(42, 5).Value = '  -2.38%  '

$ws.Cells.Item(43, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(43, 4).Value = '0.0₆0317'
$ws.Cells.Item(43, 5).Value = '  -6.33%  '

$ws.Cells.Item(44, 2).Value = 'USDe'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(44, 4).Value = '1.00'
$ws.Cells.Item(44, 5).Value = '  +0.04%  '

$ws.Cells.Item(45, 4).Value = '157.25'
$ws.Cells.Item(45, 5).Value = '  +0.20%  '

$ws.Cells.Item(46, 4).Value = '3.79'
$ws.Cells.Item(46, 5).Value = '  +0.33%  '

$ws.Cells.Item(47, 4).Value = '21.92'
$ws.Cells.Item(47, 5).Value = '  -0.32%  '

$ws.Cells.Item(48, 4).Value = '1.69'
$ws.Cells.Item(48, 5).Value = '  -2.19%  '

$ws.Cells.Item(49, 4).Value = '0.0772'
$ws.Cells.Item(49, 5).Value = '  -2.16%  '

$ws.Cells.Item(50, 2).Value = 'ARBITRUM'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(50, 4).Value = '0.574'
$ws.Cells.Item(50, 5).Value = '  +0.83%  '

$ws.Cells.Item(51, 2).Value = 'Mantle'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51, 4).Value = '0.615'
$ws.Cells.Item(51, 5).Value = '  -0.27%  '
